# The "MetKard[]" (repeated-value) columns for testKwantWrd and
# testStringField under testComplexTypeMetKard[] can't be written as a
# single CSV/XLSX cell (nested list), so those two whole columns are
# removed from the export: column M (testComplexTypeMetKard[].testStringFieldMetKard[])
# and column K (testComplexTypeMetKard[].testKwantWrdMetKard[]).
# Columns to the right of each deleted column shift left automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the rightmost column first so the left one's address ("K1")
# still refers to the intended column when it is deleted next.
$ws.Range("M1").EntireColumn.Delete()
$ws.Range("K1").EntireColumn.Delete()

# Leave the selection where Excel would after deleting those columns.
$ws.Range("M2").Select() | Out-Null
